$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update a couple of existing data rows on FormParameterData (sheet5)
#    and FormTypedata (sheet6) - renamed automation test records.
# ---------------------------------------------------------------------------
$wsParam = $wb.Worksheets.Item("FormParameterData")
$wsParam.Range("A2").Value = "Auto_561"
$wsParam.Range("A3").Value = "Auto_562"
$wsParam.Range("B2").Value = "Raf_label561"
$wsParam.Range("B3").Value = "Raf_labeL562"

$wsType = $wb.Worksheets.Item("FormTypedata")
$wsType.Range("A2").Value = "Auto-561"
$wsType.Range("A3").Value = "Auto-562"
$wsType.Range("B2").Value = "Automation561"
$wsType.Range("B3").Value = "Automation562"

# ---------------------------------------------------------------------------
# 2. Add three new reference-data sheets at the end of the workbook:
#    Department, SubDepartment and Manufacture.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsDept = $wb.Worksheets.Add($null, $lastSheet)
$wsDept.Name = "Department"
$wsDept.Columns.Item(1).ColumnWidth = 14.33
$wsDept.Columns.Item(2).ColumnWidth = 20.17
$wsDept.Columns.Item(3).ColumnWidth = 20.5

$wsDept.Range("A1").Value = "Department Id"
$wsDept.Range("B1").Value = "Description"
$wsDept.Range("C1").Value = "Metadata"
$wsDept.Range("A3").Value = "Auto-124"
$wsDept.Range("B3").Value = "Automation"
$wsDept.Range("C3").Value = "Automation"
$wsDept.Range("A2").Value = "Auto-125"
$wsDept.Range("B2").Value = "Automation"
$wsDept.Range("C2").Value = "Automation"

$wsSubDept = $wb.Worksheets.Add($null, $wsDept)
$wsSubDept.Name = "SubDepartment"
$wsSubDept.Columns.Item(1).ColumnWidth = 19.33
$wsSubDept.Columns.Item(2).ColumnWidth = 15.67

$wsSubDept.Range("A1").Value = "sub Department Id"
$wsSubDept.Range("B1").Value = "Description"
$wsSubDept.Range("C1").Value = "Metadata"

$wsManu = $wb.Worksheets.Add($null, $wsSubDept)
$wsManu.Name = "Manufacture"
$wsManu.Columns.Item(1).ColumnWidth = 16.33
$wsManu.Columns.Item(2).ColumnWidth = 20.17
$wsManu.Columns.Item(3).ColumnWidth = 19.5

$wsManu.Range("A1").Value = "Manufacture Id"
$wsManu.Range("B1").Value = "Description"
$wsManu.Range("C1").Value = "Metadata"
$wsManu.Range("C17").Value = "S"

$wsSubDept.Range("A2").Value = "Auto-SD130"
$wsSubDept.Range("B2").Value = "Automation"
$wsSubDept.Range("C2").Value = "Automation"
$wsSubDept.Range("A3").Value = "Auto-SD131"
$wsSubDept.Range("B3").Value = "Automation"
$wsSubDept.Range("C3").Value = "Automation"

$wsManu.Range("A2").Value = "Auto-MID132"
$wsManu.Range("B2").Value = "Automation"
$wsManu.Range("C2").Value = "Automation"
$wsManu.Range("A3").Value = "Auto-MID133"
$wsManu.Range("B3").Value = "Automation"
$wsManu.Range("C3").Value = "Automation"

# ---------------------------------------------------------------------------
# 3. Restore per-sheet cursor / selection positions exactly as left by the
#    author, finishing on the Manufacture sheet (which becomes the active
#    tab of the workbook).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Range("F27").Select() | Out-Null
$wsParam.Range("D19").Select() | Out-Null
$wsType.Range("C10").Select() | Out-Null
$wsDept.Range("C13").Select() | Out-Null
$wsSubDept.Range("C10").Select() | Out-Null
$wsManu.Range("B10").Select() | Out-Null
